$d = $word.ActiveDocument

# Work from the bottom of the document upward so paragraph indices for
# content above the current edit point remain stable while we edit.

# --- Paragraph 17: "This covers the core validation criteria..." ---
# Replace its text, then append two new paragraphs (blank + new item 7).
$p17 = $d.Paragraphs.Item(17)
$p17.Range.Text = "6. The system shall print out a message indicating if the email, name and password are valid or not along with an appropriate message."
$p17.Range.InsertParagraphAfter()
$p17.Range.InsertParagraphAfter()
$p19 = $d.Paragraphs.Item(19)
$p19.Range.Text = "7. The system shall call the validation functions for email, name and password and print the results."

# --- Paragraph 15: "It runs tests on various password examples..." ---
$p15 = $d.Paragraphs.Item(15)
$p15.Range.Text = '5. The system shall validate if a password contains at least one special character from the set !@#$%^&*(),.?":{}|<>`.'

# --- Paragraph 13: "It prints custom failure messages..." ---
$p13 = $d.Paragraphs.Item(13)
$p13.Range.Text = "4. The system shall validate if a password contains at least one numeric digit."

# --- Paragraphs 9-11 ("1. Checking the length...", "2. Using a regex...digit",
#     "3. Using a regex...special character") collapse into a single paragraph.
#     Delete the two trailing paragraphs (10 and, after that shifts the
#     collection, what was 11 is now also index 10), then retext paragraph 9. ---
$p10 = $d.Paragraphs.Item(10)
$p10.Range.Delete()
$p10b = $d.Paragraphs.Item(10)
$p10b.Range.Delete()
$p9 = $d.Paragraphs.Item(9)
$p9.Range.Text = "3. The system shall validate if a password is at least 8 characters long."

# --- Paragraph 7: "The test code implements automated checks..." ---
$p7 = $d.Paragraphs.Item(7)
$p7.Range.Text = "2. The system shall validate if a name is at least 8 characters long. "

# --- Paragraphs 3-5 ("- The password must be at least 8 characters long",
#     "- The password must contain at least one digit (0-9)  ",
#     "- The password must contain at least one special character...") collapse
#     into a single paragraph. Delete the two trailing paragraphs (4, then the
#     new 4), then retext paragraph 3. ---
$p4 = $d.Paragraphs.Item(4)
$p4.Range.Delete()
$p4b = $d.Paragraphs.Item(4)
$p4b.Range.Delete()
$p3 = $d.Paragraphs.Item(3)
$p3.Range.Text = "1. The system shall validate if an email address is valid."

# --- Paragraph 1: "Here are the password requirements in plain English:" ---
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Text = "Here are the natural language requirements based on the Python test code:"
